# "add genders in batumi" - extend the hotels & restaurants table on
# Sheet1 with a new "2023" column (S), mirroring the existing "2022"
# column (R): same per-row formatting, one new data point per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Write the new column's values (year header + 10 data rows + the
#        "magnitude nil / not available" marker on the totals row). ---
$ws.Range("S3").Value  = 2023
$ws.Range("S4").Value  = 58.6
$ws.Range("S5").Value  = 58.6
$ws.Range("S6").Value  = 1294
$ws.Range("S7").Value  = 1057
$ws.Range("S8").Value  = 976.4
$ws.Range("S10").Value = 12.4
$ws.Range("S11").Value = 23.4
$ws.Range("S12").Value = 2
$ws.Range("S14").Value = "_"

# Rows 9 & 13 ("Average monthly remuneration" / "Total purchases") carry
# the same floating-point representation noise as the source data for
# those rows elsewhere in the table (e.g. the 2018 column), so route them
# through a formula evaluation before freezing back to a literal value.
$ws.Range("S9").Formula = "=35.2"
$ws.Range("S9").Copy() | Out-Null
$ws.Range("S9").PasteSpecial(-4163) | Out-Null
$ws.Range("S13").Formula = "=35.2"
$ws.Range("S13").Copy() | Out-Null
$ws.Range("S13").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# --- 2. Pick up the formatting (number format / font / borders /
#        alignment) of the previous "2022" column (R) for every one of
#        those rows, so the new column visually matches the table. ---
for ($r = 3; $r -le 14; $r++) {
    $ws.Range("R$r").Copy() | Out-Null
    $ws.Range("S$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- 3. Match the selection left behind by the edit. ---
$ws.Range("S3:S14").Select()
